$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, bordered, centered) from H1 onto the new
# header cells I1:J1 before writing their text, so the new cells end up
# with the same cell style as the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I (I0) and J (IF) for rows 2-15
$data = @(
    @(3, 4),
    @(5, 6),
    @(6, 7),
    @(7, 8),
    @(4, 6),
    @(8, 9),
    @(8, 9),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(9, 9),
    @(5, 5),
    @(6, 6),
    @(3, 3)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
